$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update company_name code, and refreshed financial metrics ---
# B2 holds a numeric-looking code ("3") but must stay a text value (matches original "5" text cell),
# so force text formatting, assign, then clear the format so no stray style index is left behind.
$b2 = $ws.Range("B2")
$b2.NumberFormat = "@"
$b2.Value = "3"
$b2.ClearFormats()

# Row 2 numeric metric updates
$ws.Range("D2").Value = 0.114
$ws.Range("E2").Value = -0.207
$ws.Range("G2").Value = 0.611406844106464
$ws.Range("H2").Value = 0.611406844106464
$ws.Range("I2").Value = 0.5627376425855514
$ws.Range("J2").Value = 0.2848242278700554
$ws.Range("K2").Value = 4.887
$ws.Range("L2").Value = 0.03716349809885931
$ws.Range("M2").Value = 3.63
$ws.Range("N2").Value = 0.04514925373134328
$ws.Range("O2").Value = 0.7427869858809085
$ws.Range("P2").Value = 3.63
$ws.Range("Q2").Value = 0.04514925373134328
$ws.Range("R2").Value = 0.7427869858809085
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0
$ws.Range("U2").Value = 28.12
$ws.Range("V2").Value = 0.3497512437810945
$ws.Range("W2").Value = 0.01823529411764706
$ws.Range("X2").Value = 0.08614028408598519
$ws.Range("Y2").Value = -0.06790498996833813
$ws.Range("Z2").Value = 0.4322529748208533
$ws.Range("AA2").Value = 0.03887221374797028
$ws.Range("AB2").Value = 0.05923218467750604
$ws.Range("AC2").Value = -0.02035997092953576
$ws.Range("AD2").Value = 257.2
$ws.Range("AF2").Value = 257.2
$ws.Range("AG2").Value = 229.08
$ws.Range("AH2").Value = 0.7618483412322274
$ws.Range("AI2").Value = 0.8065224208215742
$ws.Range("AJ2").Value = 0.7402093834819696
$ws.Range("AK2").Value = 0.7878120916156545
$ws.Range("AL2").Value = 58.8
$ws.Range("AM2").Value = 58.58499999999999
$ws.Range("AN2").Value = 3.239294710327456
$ws.Range("AO2").Value = 1.258503401360544
$ws.Range("AP2").Value = 2.885138539042821
$ws.Range("AQ2").Value = 1.263121959545959

# Row 3 numeric metric updates
$ws.Range("B3").Value = 'Tunisie Leasing & Factoring Société anonyme (BVMT:TLS)'
$ws.Range("D3").Value = 0.212
$ws.Range("E3").Value = -0.0524
$ws.Range("G3").Value = 0.6057416267942584
$ws.Range("H3").Value = 0.6057416267942584
$ws.Range("I3").Value = 0.5569377990430622
$ws.Range("J3").Value = 0.2887282800302191
$ws.Range("K3").Value = 4.02
$ws.Range("L3").Value = 0.0384688995215311
$ws.Range("M3").Value = 2.27
$ws.Range("N3").Value = 0.05509708737864077
$ws.Range("O3").Value = 0.5646766169154229
$ws.Range("P3").Value = 2.27
$ws.Range("Q3").Value = 0.05509708737864077
$ws.Range("R3").Value = 0.5646766169154229
$ws.Range("S3").Value = 0
$ws.Range("T3").Value = 0
$ws.Range("X3").Value = 0.02580261134732567
$ws.Range("AB3").Value = 0.02580261134732567
$ws.Range("AL3").Value = 45.8
$ws.Range("AM3").Value = 45.8
$ws.Range("AO3").Value = 1.270742358078603
$ws.Range("AQ3").Value = 1.270742358078603

# Row 4 numeric metric updates
$ws.Range("B4").Value = 'Arab Tunisian Lease S.A. (BVMT:ATL)'
$ws.Range("D4").Value = 0.114
$ws.Range("E4").Value = -0.291
$ws.Range("G4").Value = 0.7339055793991417
$ws.Range("H4").Value = 0.7339055793991417
$ws.Range("I4").Value = 0.6781115879828327
$ws.Range("J4").Value = 0.3390557939914163
$ws.Range("K4").Value = 0.402
$ws.Range("L4").Value = 0.01725321888412017
$ws.Range("O4").Value = -0
$ws.Range("R4").Value = -0
$ws.Range("U4").Value = 4.82
$ws.Range("V4").Value = 0.2510416666666667
$ws.Range("W4").Value = 0.01435714285714286
$ws.Range("X4").Value = 0.1287505422913847
$ws.Range("Y4").Value = -0.1143933994342418
$ws.Range("Z4").Value = 0.1146484278895832
$ws.Range("AA4").Value = 0.03887221374797028
$ws.Range("AB4").Value = 0.05923218467750604
$ws.Range("AC4").Value = -0.02035997092953576
$ws.Range("AD4").Value = 159.7
$ws.Range("AF4").Value = 159.7
$ws.Range("AG4").Value = 154.88
$ws.Range("AH4").Value = 0.8926774734488542
$ws.Range("AI4").Value = 0.8244708311822405
$ws.Range("AJ4").Value = 0.8897058823529412
$ws.Range("AK4").Value = 0.81999152901313
$ws.Range("AL4").Value = 13
$ws.Range("AM4").Value = 12.785
$ws.Range("AN4").Value = 9.919254658385091
$ws.Range("AO4").Value = 1.215384615384616
$ws.Range("AP4").Value = 9.619875776397514
$ws.Range("AQ4").Value = 1.235823230348064

# Row 5 numeric metric updates
$ws.Range("B5").Value = 'Best Lease SA (BVMT:BL)'
$ws.Range("D5").Value = 0.0298
$ws.Range("E5").Value = -0.207
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0.465
$ws.Range("L5").Value = 0.1256756756756757
$ws.Range("M5").Value = 1.36
$ws.Range("N5").Value = 0.068
$ws.Range("O5").Value = 2.924731182795699
$ws.Range("P5").Value = 1.36
$ws.Range("Q5").Value = 0.068
$ws.Range("R5").Value = 2.924731182795699
$ws.Range("S5").Value = 0
$ws.Range("T5").Value = 0
$ws.Range("U5").Value = 23.3
$ws.Range("V5").Value = 1.165
$ws.Range("W5").Value = 0.01823529411764706
$ws.Range("X5").Value = 0.08614028408598519
$ws.Range("Y5").Value = -0.06790498996833813
$ws.Range("Z5").Value = 0.03663729082087335
$ws.Range("AA5").Value = 0
$ws.Range("AB5").Value = 0.06358585940312199
$ws.Range("AC5").Value = -0.06358585940312199
$ws.Range("AD5").Value = 97.5
$ws.Range("AF5").Value = 97.5
$ws.Range("AG5").Value = 74.2
$ws.Range("AH5").Value = 0.8297872340425532
$ws.Range("AI5").Value = 0.7787539936102236
$ws.Range("AJ5").Value = 0.7876857749469215
$ws.Range("AK5").Value = 0.7281648675171737
$ws.Range("AL5").Value = 0
$ws.Range("AM5").Value = 0

# Row 5: these trailing metric columns no longer apply for the new company in this row
$ws.Range("AN5").ClearContents()
$ws.Range("AO5").ClearContents()
$ws.Range("AP5").ClearContents()
$ws.Range("AQ5").ClearContents()

# Consolidate table: rows 6 and 7 companies were merged/reassigned into rows 2-5 above, so remove them
$ws.Rows("6:7").Delete()
